# Applies the "almost all programs are ready" edit to slide 2:
#  - "順時鐘轉" -> "向下扳動"            (shape 6, 文字方塊 1)
#  - "逆時鐘轉" -> "向上扳動"            (shape 7, 文字方塊 7)
#  - "按旋鈕"   -> "向右扳動" + widen box (shape 8, 文字方塊 8)
#  - "按按鈕"   -> "向左扳動" + widen box (shape 9, 文字方塊 9)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Shape 6: 文字方塊 1 (順時鐘轉 -> 向下扳動)
$s.Shapes.Item(6).TextFrame.TextRange.Text = "向下扳動"

# Shape 7: 文字方塊 7 (逆時鐘轉 -> 向上扳動)
$s.Shapes.Item(7).TextFrame.TextRange.Text = "向上扳動"

# Shape 8: 文字方塊 8 (按旋鈕 -> 向右扳動), widened from 877163 EMU to 1107996 EMU
$s.Shapes.Item(8).TextFrame.TextRange.Text = "向右扳動"
$s.Shapes.Item(8).Width = 1107996 / 12700

# Shape 9: 文字方塊 9 (按按鈕 -> 向左扳動), widened from 877163 EMU to 1107996 EMU
$s.Shapes.Item(9).TextFrame.TextRange.Text = "向左扳動"
$s.Shapes.Item(9).Width = 1107996 / 12700
